# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Kujata_Profits workbook (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 199.66667
$ws.Range("I8").Value = 199.66667
$ws.Range("K8").Value = 599.00001
$ws.Range("M8").Value = -460.00001
$ws.Range("H18").Value = 2271.2856
$ws.Range("I18").Value = 2699.75
$ws.Range("J18").Value = 1700
$ws.Range("K18").Value = 2699.75
$ws.Range("L18").Value = 1700
$ws.Range("M18").Value = -2415.75
$ws.Range("N18").Value = -2268
$ws.Range("H40").Value = 2666.1667
$ws.Range("I40").Value = 2199.4
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 2199.4
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -2024.4
$ws.Range("N40").Value = -5350
$ws.Range("H82").Value = 1000
$ws.Range("I82").Value = 1000
$ws.Range("K82").Value = 3000
$ws.Range("M82").Value = -2594
$ws.Range("H85").Value = 1000
$ws.Range("I85").Value = 1000
$ws.Range("K85").Value = 3000
$ws.Range("M85").Value = -1596
$ws.Range("H98").Value = 4512.8887
$ws.Range("I98").Value = 4602.6
$ws.Range("J98").Value = 4400.75
$ws.Range("K98").Value = 4602.6
$ws.Range("L98").Value = 4400.75
$ws.Range("M98").Value = -3104.6
$ws.Range("N98").Value = -7396.75
$ws.Range("H112").Value = 2755.4736
$ws.Range("J112").Value = 2755.4736
$ws.Range("L112").Value = 8266.4208
$ws.Range("N112").Value = -10482.4208
$ws.Range("H122").Value = 4512.8887
$ws.Range("I122").Value = 4602.6
$ws.Range("J122").Value = 4400.75
$ws.Range("K122").Value = 13807.8
$ws.Range("L122").Value = 13202.25
$ws.Range("M122").Value = -11357.8
$ws.Range("N122").Value = -18102.25
$ws.Range("H125").Value = 1049.75
$ws.Range("I125").Value = 1049.75
$ws.Range("K125").Value = 9447.75
$ws.Range("M125").Value = -6987.75
$ws.Range("H135").Value = 793.5333000000001
$ws.Range("J135").Value = 4000
$ws.Range("L135").Value = 36000
$ws.Range("N135").Value = -41070
$ws.Range("H136").Value = 72890
$ws.Range("J136").Value = 72890
$ws.Range("L136").Value = 72890
$ws.Range("N136").Value = -83090
$ws.Range("H137").Value = 1415.5834
$ws.Range("I137").Value = 1317.6666
$ws.Range("K137").Value = 3952.9998
$ws.Range("M137").Value = -1402.9998
$ws.Range("H141").Value = 865.625
$ws.Range("I141").Value = 865.625
$ws.Range("K141").Value = 2596.875
$ws.Range("M141").Value = 2583.125

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3525.4807
$ws.Range("I32").Value = 3112.383
$ws.Range("J32").Value = 7408.6
$ws.Range("K32").Value = 3112.383
$ws.Range("L32").Value = 7408.6
$ws.Range("M32").Value = -2825.383
$ws.Range("N32").Value = -7982.6
$ws.Range("H61").Value = 1285.0714
$ws.Range("I61").Value = 999.3333
$ws.Range("K61").Value = 999.3333
$ws.Range("M61").Value = -787.3333
$ws.Range("H74").Value = 1204.7693
$ws.Range("I74").Value = 531
$ws.Range("K74").Value = 531
$ws.Range("M74").Value = 343
$ws.Range("H77").Value = 1204.7693
$ws.Range("I77").Value = 531
$ws.Range("K77").Value = 2655
$ws.Range("M77").Value = 1713
$ws.Range("H109").Value = 6000
$ws.Range("J109").Value = 6000
$ws.Range("L109").Value = 6000
$ws.Range("N109").Value = -8774
$ws.Range("H132").Value = 2751
$ws.Range("I132").Value = 2665.7273
$ws.Range("J132").Value = 2823.1538
$ws.Range("K132").Value = 7997.1819
$ws.Range("L132").Value = 8469.4614
$ws.Range("M132").Value = -5467.1819
$ws.Range("N132").Value = -13529.4614
$ws.Range("H136").Value = 1285.0714
$ws.Range("I136").Value = 999.3333
$ws.Range("K136").Value = 2997.9999
$ws.Range("M136").Value = -447.9998999999998

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2313.8572
$ws.Range("I20").Value = 2312.5
$ws.Range("J20").Value = 2315.6667
$ws.Range("K20").Value = 2312.5
$ws.Range("L20").Value = 2315.6667
$ws.Range("M20").Value = -2065.5
$ws.Range("N20").Value = -2809.6667
$ws.Range("H107").Value = 1877.3077
$ws.Range("I107").Value = 1399.5
$ws.Range("K107").Value = 1399.5
$ws.Range("M107").Value = 520.5
$ws.Range("H132").Value = 64666.332
$ws.Range("J132").Value = 64666.332
$ws.Range("L132").Value = 64666.332
$ws.Range("N132").Value = -74786.33199999999
$ws.Range("H134").Value = 9424.134
$ws.Range("I134").Value = 1036.3
$ws.Range("J134").Value = 26199.8
$ws.Range("K134").Value = 3108.9
$ws.Range("L134").Value = 78599.39999999999
$ws.Range("M134").Value = -573.8999999999996
$ws.Range("N134").Value = -83669.39999999999

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 248.8
$ws.Range("I7").Value = 173.5
$ws.Range("J7").Value = 550
$ws.Range("K7").Value = 173.5
$ws.Range("L7").Value = 550
$ws.Range("M7").Value = -60.5
$ws.Range("N7").Value = -776
$ws.Range("H22").Value = 477.77777
$ws.Range("I22").Value = 383.33334
$ws.Range("J22").Value = 666.6667
$ws.Range("K22").Value = 383.33334
$ws.Range("L22").Value = 666.6667
$ws.Range("M22").Value = -33.33334000000002
$ws.Range("N22").Value = -1366.6667
$ws.Range("H58").Value = 1288.0834
$ws.Range("I58").Value = 1009.1111
$ws.Range("K58").Value = 1009.1111
$ws.Range("M58").Value = -806.1111
$ws.Range("H132").Value = 5628.6665
$ws.Range("I132").Value = 7534.2354
$ws.Range("J132").Value = 3136.7693
$ws.Range("K132").Value = 22602.7062
$ws.Range("L132").Value = 9410.3079
$ws.Range("M132").Value = -20072.7062
$ws.Range("N132").Value = -14470.3079
$ws.Range("H136").Value = 1288.0834
$ws.Range("I136").Value = 1009.1111
$ws.Range("K136").Value = 3027.3333
$ws.Range("M136").Value = -477.3332999999998

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1335.6818
$ws.Range("I5").Value = 1335.6818
$ws.Range("K5").Value = 4007.0454
$ws.Range("M5").Value = -3895.0454
$ws.Range("H103").Value = 7117.7646
$ws.Range("J103").Value = 10717
$ws.Range("L103").Value = 32151
$ws.Range("N103").Value = -33909
$ws.Range("H131").Value = 12821765
$ws.Range("I131").Value = 250000290
$ws.Range("J131").Value = 1304.9324
$ws.Range("K131").Value = 750000870
$ws.Range("L131").Value = 3914.7972
$ws.Range("M131").Value = -749995830
$ws.Range("N131").Value = -13994.7972
$ws.Range("H135").Value = 1335.6818
$ws.Range("I135").Value = 1335.6818
$ws.Range("K135").Value = 12021.1362
$ws.Range("M135").Value = -9486.136200000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()
$ws.Range("H70").Value = 37506370
$ws.Range("I70").Value = 35720060
$ws.Range("K70").Value = 35720060
$ws.Range("M70").Value = -35719790
$ws.Range("H73").Value = 37506370
$ws.Range("I73").Value = 35720060
$ws.Range("K73").Value = 35720060
$ws.Range("M73").Value = -35719124
$ws.Range("H122").Value = 3199.6667
$ws.Range("I122").Value = 3199.6667
$ws.Range("K122").Value = 9599.000100000001
$ws.Range("M122").Value = -7149.000100000001
$ws.Range("H132").Value = 2239.6562
$ws.Range("I132").Value = 1850.125
$ws.Range("J132").Value = 3408.25
$ws.Range("K132").Value = 5550.375
$ws.Range("L132").Value = 10224.75
$ws.Range("M132").Value = -3020.375
$ws.Range("N132").Value = -15284.75

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2216.75
$ws.Range("I7").Value = 2033.4286
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 2033.4286
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -1921.4286
$ws.Range("N7").Value = -3724
$ws.Range("H40").Value = 2802
$ws.Range("I40").Value = 2817.8333
$ws.Range("J40").Value = 2754.5
$ws.Range("K40").Value = 2817.8333
$ws.Range("L40").Value = 2754.5
$ws.Range("M40").Value = -2681.8333
$ws.Range("N40").Value = -3026.5
$ws.Range("H116").Value = 25680
$ws.Range("J116").Value = 25680
$ws.Range("L116").Value = 25680
$ws.Range("N116").Value = -34858
$ws.Range("H126").Value = 2216.75
$ws.Range("I126").Value = 2033.4286
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 6100.2858
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -3630.2858
$ws.Range("N126").Value = -15440
$ws.Range("H132").Value = 21989.918
$ws.Range("I132").Value = 1335.5
$ws.Range("J132").Value = 49529.145
$ws.Range("K132").Value = 4006.5
$ws.Range("L132").Value = 148587.435
$ws.Range("M132").Value = -1476.5
$ws.Range("N132").Value = -153647.435

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 352002.5
$ws.Range("I14").Value = 352002.5
$ws.Range("K14").Value = 352002.5
$ws.Range("M14").Value = -351834.5
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H122").Value = 13687133
$ws.Range("I122").Value = 16253158
$ws.Range("J122").Value = 1666.6666
$ws.Range("K122").Value = 48759474
$ws.Range("L122").Value = 4999.9998
$ws.Range("M122").Value = -48757024
$ws.Range("N122").Value = -9899.9998
$ws.Range("H132").Value = 2731.2964
$ws.Range("I132").Value = 2152.4092
$ws.Range("J132").Value = 5278.4
$ws.Range("K132").Value = 6457.2276
$ws.Range("L132").Value = 15835.2
$ws.Range("M132").Value = -3927.2276
$ws.Range("N132").Value = -20895.2
$ws.Range("H136").Value = 1207.7826
$ws.Range("I136").Value = 787.1177
$ws.Range("K136").Value = 2361.3531
$ws.Range("M136").Value = 188.6468999999997
$ws.Range("H138").Value = 34986
$ws.Range("J138").Value = 34986
$ws.Range("L138").Value = 34986
$ws.Range("N138").Value = -45266
